$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all previous data rows (A2:A36) first
$ws.Range("A2:A36").ClearContents()

# New consolidated card data (tuples as strings), rows 2 through 10
$ws.Range("A2").Value = "('Castigate', ['{W}{B}', 'Sorcery', 'Target opponent reveals their hand. You choose a nonland card from it and exile that card.'])"
$ws.Range("A3").Value = "('Coiling Oracle', ['{G}{U}', 'Creature " + [char]0x2014 + " Snake Elf Druid', 'When Coiling Oracle enters the battlefield, reveal the top card of your library. If it" + [char]0x2019 + "s a land card, put it onto the battlefield. Otherwise, put that card into your hand.', '1/1'])"
$ws.Range("A4").Value = "('Forest', ['Basic Land " + [char]0x2014 + " Forest', '({T}: Add {G}.)'])"
$ws.Range("A5").Value = "('Island', ['Basic Land " + [char]0x2014 + " Island', '({T}: Add {U}.)'])"
$ws.Range("A6").Value = "('Mountain', ['Basic Land " + [char]0x2014 + " Mountain', '({T}: Add {R}.)'])"
$ws.Range("A7").Value = "('Plains', ['Basic Land " + [char]0x2014 + " Plains', '({T}: Add {W}.)'])"
$ws.Range("A8").Value = "('Surging Flame', ['{1}{R}', 'Instant', 'Ripple 4 (When you cast this spell, you may reveal the top four cards of your library. You may cast spells with the same name as this spell from among those cards without paying their mana costs. Put the rest on the bottom of your library.)', 'Surging Flame deals 2 damage to any target.'])"
$ws.Range("A9").Value = "('Swamp', ['Basic Land " + [char]0x2014 + " Swamp', '({T}: Add {B}.)'])"
$ws.Range("A10").Value = "('Wee Dragonauts', ['{1}{U}{R}', 'Creature " + [char]0x2014 + " Faerie Wizard', 'Flying', 'Whenever you cast an instant or sorcery spell, Wee Dragonauts gets +2/+0 until end of turn.', '1/3'])"

# Remove the now-unused rows 11-36 so the used range shrinks to A1:A10
$ws.Range("A11:A36").Delete()
